$d = $word.ActiveDocument

# The second paragraph holds a Word field: { m:self. } (displayed via
# fldChar begin/instrText/fldChar end). The parser was updated to use
# TokenIteratorFieldRewriterSplit, so the field is rewritten as plain
# literal text runs "{", "m", ":", "self" (kept colored) and ".}".

$fld = $d.Fields.Item(1)
$p2 = $d.Paragraphs.Item(2)

# Drop the field itself (removes fldChar/instrText runs), leaving an
# empty paragraph in its place.
$fld.Delete()

$r = $p2.Range
$r.Collapse(1)
$r.Text = "{m:self.}"

# Re-split "self" into its own colored run, matching the
# accent6/BF-shade orange used by the original field code run.
$start = $p2.Range.Start
$selfRange = $d.Range($start + 3, $start + 7)
$selfRange.Font.TextColor.ObjectThemeColor = 9
